# Update profit worksheet after running on 2025-09-30:
# append a new row (44) with the date and profit figure for that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A in this sheet stores dates as literal text (e.g. "09/29/2025"
# in row 43), not as real Excel date serials. Force the cell to Text
# format before assigning so Excel doesn't auto-convert the string into
# a date value, then drop back to the default "Normal" style so no new
# style gets stamped onto the cell (matching the unstyled cells above it).
$dateCell = $ws.Range("A44")
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/30/2025"
$dateCell.Style = "Normal"

$ws.Range("B44").Value = 14397.21
